$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 98, column A: refreshed timestamp with a tiny floating-point
# precision difference versus the original value.
$ws.Range("A98").Value = 44411.76825210301

# New row 99: freshly retrieved data row.
$ws.Range("A99").Value = 44412.77001802667
$ws.Range("B99").Value = 79702
$ws.Range("C99").Value = 67607
$ws.Range("D99").Value = 3596
$ws.Range("E99").Value = 2287
$ws.Range("F99").Value = 1645
$ws.Range("G99").Value = 21041
$ws.Range("H99").Value = 1615
$ws.Range("I99").Value = 927
$ws.Range("J99").Value = 200
